$d = $word.ActiveDocument

# --- Helpers -------------------------------------------------------------
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-BodyXml($range, $bodyFragment) {
    $full = $pkgHeader + '<w:body>' + $bodyFragment + '</w:body>' + $pkgFooter
    $range.InsertXML($full)
}

# --- Locate the "Query" paragraph and the blank paragraph right after it -
$queryParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Query`r") {
        $queryParaIndex = $i
        break
    }
}
if ($queryParaIndex -eq -1) {
    throw "Could not find the standalone 'Query' paragraph"
}

$pQuery = $d.Paragraphs.Item($queryParaIndex)
$pAfterQuery = $d.Paragraphs.Item($queryParaIndex + 1)
$queryRange = $d.Range($pQuery.Range.Start, $pAfterQuery.Range.End)

$queryFragment = '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Query</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'

Insert-BodyXml $queryRange $queryFragment

# --- Append the new Overpass-query paragraphs at the very end of the body
$endRange = $d.Range($d.Content.End, $d.Content.End)

$tailFragment = '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>out:csv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>length,value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>)];</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>area[name="Nordrhein-Westfalen"]-&gt;.a;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>area[name="Bielefeld"]-&gt;.b;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>way[bicycle](</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>area.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>)(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>area.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>for (t["bicycle"])</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>{</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">  make stat value=_.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>val,length</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>=sum(length());</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">  out;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}</w:t></w:r></w:p>'

Insert-BodyXml $endRange $tailFragment

Write-Output "done"
